$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "24.654.52"
$ws.Cells.Item(2, 5).Value = "  -0.31%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.696.95"
$ws.Cells.Item(3, 5).Value = "  +0.10%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.008"
$ws.Cells.Item(4, 5).Value = "  +0.75%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "315.65"
$ws.Cells.Item(5, 5).Value = "  -0.25%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.012"
$ws.Cells.Item(6, 5).Value = "  +1.14%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3929"
$ws.Cells.Item(7, 5).Value = "  -0.46%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.4063"
$ws.Cells.Item(8, 5).Value = "  +0.17%  "

$ws.Cells.Item(9, 5).Value = "  +1.60%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.008"
$ws.Cells.Item(10, 5).Value = "  +0.75%  "

$ws.Cells.Item(11, 5).Value = "  -0.43%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08786"
$ws.Cells.Item(12, 5).Value = "  -0.94%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "7.654"
$ws.Cells.Item(13, 5).Value = "  +5.87%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "24.42"
$ws.Cells.Item(14, 5).Value = "  +3.33%  "

$ws.Cells.Item(15, 5).Value = "  +3.19%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "8.003"
$ws.Cells.Item(16, 5).Value = "  -0.57%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.712.85"
$ws.Cells.Item(17, 5).Value = "  +1.15%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "98.74"
$ws.Cells.Item(18, 5).Value = "  -1.25%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.07120"
$ws.Cells.Item(19, 5).Value = "  +1.56%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "19.83"
$ws.Cells.Item(20, 5).Value = "  +1.01%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "7.402"

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.011"
$ws.Cells.Item(22, 5).Value = "  +1.00%  "

$ws.Cells.Item(23, 5).Value = "  -0.24%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "24.656.23"
$ws.Cells.Item(24, 5).Value = "  -0.30%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.052"
$ws.Cells.Item(25, 5).Value = "  -6.45%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.363"
$ws.Cells.Item(26, 5).Value = "  +0.44%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "22.74"
$ws.Cells.Item(27, 5).Value = "  +0.00%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "163.81"
$ws.Cells.Item(28, 5).Value = "  +0.34%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.473"
$ws.Cells.Item(29, 5).Value = "  +13.20%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "137.97"
$ws.Cells.Item(30, 5).Value = "  +1.31%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "5.249"
$ws.Cells.Item(31, 5).Value = "  +1.01%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.898.05"
$ws.Cells.Item(32, 5).Value = "  +0.98%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.08839"
$ws.Cells.Item(33, 5).Value = "  +2.88%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "7.484"
$ws.Cells.Item(34, 5).Value = "  +4.72%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.053"
$ws.Cells.Item(35, 5).Value = "  -1.56%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.988"
$ws.Cells.Item(36, 5).Value = "  +2.97%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.02926"
$ws.Cells.Item(37, 5).Value = "  +6.90%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.2735"
$ws.Cells.Item(38, 5).Value = "  -0.48%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "10.82"
$ws.Cells.Item(39, 5).Value = "  -6.84%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "14.32"
$ws.Cells.Item(40, 5).Value = "  -1.31%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.09157"
$ws.Cells.Item(41, 5).Value = "  -0.88%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.7888"
$ws.Cells.Item(42, 5).Value = "  +2.73%  "

$ws.Cells.Item(43, 5).Value = "  -0.41%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "16.70"
$ws.Cells.Item(44, 5).Value = "  +4.76%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.7217"
$ws.Cells.Item(45, 5).Value = "  +0.20%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.566"
$ws.Cells.Item(46, 5).Value = "  -0.40%  "

$ws.Cells.Item(47, 5).Value = "  +0.07%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.006"
$ws.Cells.Item(48, 5).Value = "  +0.63%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.331"
$ws.Cells.Item(49, 5).Value = "  +0.38%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "139.39"
$ws.Cells.Item(50, 5).Value = "  -0.07%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "91.91"
$ws.Cells.Item(51, 5).Value = "  +1.97%  "
